$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 is the Schottky Barrier Diode (row 5): update its footprint and LCSC part #.
$ws.Range("C5").Value = "SOD-523"
$ws.Range("D5").Value = "C345957"

# J1,J2 4-Pin Connector (row 6): update its LCSC part #.
$ws.Range("D6").Value = "C145956"

# Reflect the author's new active selection (whole row 6 selected, active cell A6).
$ws.Rows("6:6").Select()
